# Fixed range name issue
# - Renamed sheets to clearer names (Base -> Main Sheet, Range -> Range Madness,
#   Sheet2 -> Another Sheet)
# - Added a new "Detached " worksheet right after "Main Sheet"
# - Formulas referencing the renamed sheets are auto-updated by Excel

$wb = $excel.ActiveWorkbook

$base       = $wb.Worksheets.Item("Base")
$rangeSheet = $wb.Worksheets.Item("Range")
$sheet2     = $wb.Worksheets.Item("Sheet2")

$base.Name       = "Main Sheet"
$rangeSheet.Name = "Range Madness"
$sheet2.Name     = "Another Sheet"

# Insert the new "Detached " sheet right after "Main Sheet"
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $base)
$new.Name = "Detached "

$new.Range("A1").Value = "Range"
$new.Range("A2").Value = 1
$new.Range("A3").Value = 2
$new.Range("A4").Value = 3
$new.Range("A5").Formula = "=SUM(A2:A4)"

$new.Range("A14").Value = 666
$new.Range("B14").Formula = "=A14"
$new.Range("C14").Formula = "=B14"
$new.Range("D14").Formula = "=C14"

$new.Range("A15").Select() | Out-Null

# Re-fetch the sheet objects by their final names: inserting the new
# sheet shifts worksheet positions, so stale references grabbed before
# the insert no longer point at the sheets we renamed.
$mainSheet    = $wb.Worksheets.Item("Main Sheet")
$rangeMadness = $wb.Worksheets.Item("Range Madness")

# "Range Madness" selection moved from F12 to B24
$rangeMadness.Activate() | Out-Null
$rangeMadness.Range("B24").Select() | Out-Null

# Keep "Main Sheet" as the active/selected tab with its new selection
$mainSheet.Activate() | Out-Null
$mainSheet.Range("F38").Select() | Out-Null
